$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. D24: "AENP + (EC1+EC2+EC3+EC4) - EC5" -> "Same as MEPP" ---
$ws.Range("D24").Value = "Same as MEPP"

# --- 2. D23: update the bold/green run of the rich-text comment ---
$cell = $ws.Range("D23")
$full = $cell.Value2

$oldMid = "MEPP - (EC1+EC2+EC3+EC4) + EC5" + "`n" + "(These EC values should be calculated using MEPP X % value instead of the more complex EC % formula)"
$newMid = "MEPP - (MC1+MC2+MC3+MC4) + MC5" + "`n" + "(These MC values should be calculated using EC `$ amounts or MEPP x EC % value instead of the more complex original EC % formula)"

$start = $full.IndexOf($oldMid)
$newFull = $full.Substring(0, $start) + $newMid + $full.Substring($start + $oldMid.Length)

$cell.Value2 = $newFull

# Re-apply the bold / green / Arial 12 formatting to the replaced run
$midChars = $cell.Characters($start + 1, $newMid.Length)
$midChars.Font.Name = "Arial"
$midChars.Font.Size = 12
$midChars.Font.Bold = $true
$midChars.Font.Color = 5287936

# Keep the trailing run ("Recalculate Adjusted EPP ...") on Arial 12 like before
$tailStart = $start + $newMid.Length
$tailLen = $full.Length - ($start + $oldMid.Length)
$tailChars = $cell.Characters($tailStart + 1, $tailLen)
$tailChars.Font.Name = "Arial"
$tailChars.Font.Size = 12

# --- 3. Row 23 height: 153.75 -> 177 (custom height) ---
$ws.Rows.Item(23).RowHeight = 177

# --- 4. Selection moves from D26 to E23 ---
[void]$ws.Range("E23").Select()
